# Portfolio_Analysis_Sheets.xlsx -- add Simple_Return_% column (K) to the
# Section 1 strategy-performance table and a matching "Portfolio Simple
# Return %" summary row in Section 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # Final_Portfolio_Analysis

# ---------------------------------------------------------------------
# 1) New column K: same width as the other data columns, same per-row
#    formatting as column J (header fill/border styles for rows 16-26).
# ---------------------------------------------------------------------
$ws.Columns.Item(11).ColumnWidth = 15.17

$ws.Range("J16:J26").Copy()
$ws.Range("K16:K26").PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(16, 11).Value = $null

# ---------------------------------------------------------------------
# 2) Shift the strategy data Capital/Profit/Max_DD/Years one column to
#    the right (G->H, H->I, I->J, J->K) to make room for the new
#    Simple_Return_% column at F, with Return_% moving to G.
#    Work right-to-left so no source values are clobbered.
# ---------------------------------------------------------------------
$ws.Range("J18:J25").Copy()
$ws.Range("K18:K25").PasteSpecial(-4163)   # xlPasteValues

$ws.Range("I18:I25").Copy()
$ws.Range("J18:J25").PasteSpecial(-4163)

$ws.Range("H18:H25").Copy()
$ws.Range("I18:I25").PasteSpecial(-4163)

$ws.Range("G18:G25").Copy()
$ws.Range("H18:H25").PasteSpecial(-4163)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3) Header row 17: Strategy | Pair_Method | Pairs | Sharpe | XIRR_% |
#    Simple_Return_% | Return_% | Capital | Profit | Max_DD | Years
# ---------------------------------------------------------------------
$ws.Range("F17").Value = "Simple_Return_%"
$ws.Range("G17").Value = "Return_%"
$ws.Range("H17").Value = "Capital"
$ws.Range("I17").Value = "Profit"
$ws.Range("J17").Value = "Max_DD"
$ws.Range("K17").Value = "Years"

# ---------------------------------------------------------------------
# 4) Data rows 18-25: new Simple_Return_% (F) and Return_% (G) formulas.
# ---------------------------------------------------------------------
$ws.Range("F18:F25").Formula = '=ROUND(IF(H18>0, (I18/H18/K18)*100, 0), 2)'
$ws.Range("G18:G25").Formula = '=ROUND(IF(H18>0, (I18/H18)*100, 0), 2)'

# ---------------------------------------------------------------------
# 5) TOTAL row 26: Capital/Profit/Max_DD totals now live in H/I/J;
#    G26 (old Capital total) and K26 (Years) stay blank.
# ---------------------------------------------------------------------
$ws.Range("G26").ClearContents()
$ws.Range("H26").Formula = '=ROUND(SUM(H18:H25), 2)'
$ws.Range("I26").Formula = '=ROUND(SUM(I18:I25), 2)'
$ws.Range("J26").Formula = '=ROUND(SUM(J18:J25), 2)'
$ws.Range("K26").ClearContents()

# ---------------------------------------------------------------------
# 6) Section 2 allocation table (rows 31-38): Capital column reference
#    moved from G to H.
# ---------------------------------------------------------------------
$ws.Range("C31:C38").Formula = '=ROUND((B31/100)*$H$26, 2)'
$ws.Range("D31:D38").Formula = '=ROUND(H18, 2)'

# ---------------------------------------------------------------------
# 7) Section 3 summary: Capital Required / Max Drawdown references.
# ---------------------------------------------------------------------
$ws.Range("B44").Formula = '=ROUND(H26, 2)'
$ws.Range("B45").Formula = '=ROUND(J26, 2)'

# ---------------------------------------------------------------------
# 8) Insert a new "Portfolio Simple Return %" row right after
#    "Portfolio XIRR" (row 48), pushing the rest of Section 3 down by
#    one row.
# ---------------------------------------------------------------------
$ws.Rows.Item(48).Insert()

$ws.Range("A49:C49").Copy()
$ws.Range("A48:C48").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A48").Value = "Portfolio Simple Return %"
$ws.Range("B48").Formula = '=ROUND(SUMPRODUCT(B31:B38/100, F18:F25), 2)'
$ws.Range("C48").Value = "%"

# Rows 49/50 kept their old formulas verbatim on shift; point them at
# the new Profit/Capital columns (I/H) instead of the old H/G.
$ws.Range("B49").Formula = '=ROUND(I26, 2)'
$ws.Range("B50").Formula = '=ROUND(IF(H26>0, (I26/H26)*100, 0), 2)'

# ---------------------------------------------------------------------
# 9) Re-merge the Section 1 banner across the new column.
# ---------------------------------------------------------------------
$ws.Range("A16:J16").UnMerge()
$ws.Range("A16:K16").Merge()
